# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on row 2 of the
# zh-cn and de-de sheets to reflect the new report generation times.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-13 19:12:32"
$wsZhCn.Range("H2").Value = "2016-03-13 19:12:50"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-13 19:12:35"
$wsDeDe.Range("H2").Value = "2016-03-13 19:12:56"
